$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Julio de 2020 a las 07:05"

# --- Swap the "Islas Malvinas" / "Groenlandia" rows (209-212 context: Papua Nueva Guinea, Islas Malvinas, Groenlandia, Montserrat) ---
# Row 210 currently "Islas Malvinas", row 211 currently "Groenlandia" -> swap so row 210 becomes "Groenlandia", row 211 becomes "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Pakistan (row 15): updated case counts ---
$ws.Range("B15").Value = 265083
$ws.Range("C15").Value = 1587
$ws.Range("D15").Value = 205929
$ws.Range("E15").Value = 53555
$ws.Range("G15").Value = 31
$ws.Range("H15").Value = 5599

# --- Kirguistan (row 58): updated case counts ---
$ws.Range("B58").Value = 27143
$ws.Range("C58").Value = 611
$ws.Range("D58").Value = 13109
$ws.Range("E58").Value = 12997
$ws.Range("G58").Value = 34
$ws.Range("H58").Value = 1037

# --- Tailandia (row 104): updated case counts ---
$ws.Range("B104").Value = 3250
$ws.Range("C104").Value = 1
$ws.Range("E104").Value = 96

# --- Butan (row 187): updated case counts ---
$ws.Range("B187").Value = 89
$ws.Range("C187").Value = 2
$ws.Range("D187").Value = 82
